$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted above the existing row 102,
# shifting all subsequent rows (old 102-128) down by one (new 103-129).
$ws.Rows.Item(102).Insert()

$ws.Range("A102").Value = 8
$ws.Range("B102").Value = "Terminal La Palmera de La Serena"
$ws.Range("C102").Value = "Coquimbo"
$ws.Range("D102").Value = 44855
$ws.Range("E102").Value = 4
$ws.Range("F102").Value = 100112052
$ws.Range("G102").Value = "Albahaca"
$ws.Range("H102").Value = "Sin especificar"
$ws.Range("I102").Value = "Primera"
$ws.Range("J102").Value = 1000
$ws.Range("K102").Value = 4000
$ws.Range("L102").Value = 4500
$ws.Range("M102").Value = 4250
$ws.Range("N102").Value = "$/paquete"
$ws.Range("O102").Value = "Región de Arica y Parinacota"
$ws.Range("P102").Value = 4250
$ws.Range("Q102").Value = 1
$ws.Range("R102").Value = "Hortaliza"

# Preserve the date-style formatting (matching the rest of column D) on the new cell.
$ws.Range("D102").NumberFormat = $ws.Range("D103").NumberFormat
